# Reporting Organisation Group (es) codelist: the "group-code" and
# "group-name" columns (C and D) were swapped - the ISO/IATI code now
# lives in column D and the Spanish display name in column C, for the
# header row as well as every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)

    $cVal = $cCell.Value2
    $dVal = $dCell.Value2

    $cCell.Value2 = $dVal
    $dCell.Value2 = $cVal
}
